$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rig_by_weight")

# --- Value updates (order controls shared-string insertion order: "0", "-3", "-1") ---

# "0" replaces "BASE" in C9:D11 and C31:D33
$ws.Range("C9").Value = "0"
$ws.Range("D9").Value = "0"
$ws.Range("C10").Value = "0"
$ws.Range("D10").Value = "0"
$ws.Range("C11").Value = "0"
$ws.Range("D11").Value = "0"
$ws.Range("C31").Value = "0"
$ws.Range("D31").Value = "0"
$ws.Range("C32").Value = "0"
$ws.Range("D32").Value = "0"
$ws.Range("C33").Value = "0"
$ws.Range("D33").Value = "0"

# "-3" replaces ">-2" in D2:D6 and D22:D28
$ws.Range("D2").Value = "-3"
$ws.Range("D3").Value = "-3"
$ws.Range("D4").Value = "-3"
$ws.Range("D5").Value = "-3"
$ws.Range("D6").Value = "-3"
$ws.Range("D22").Value = "-3"
$ws.Range("D23").Value = "-3"
$ws.Range("D24").Value = "-3"
$ws.Range("D25").Value = "-3"
$ws.Range("D26").Value = "-3"
$ws.Range("D27").Value = "-3"
$ws.Range("D28").Value = "-3"

# "-1" (as text) replaces numeric -1 in D7:D8
$ws.Range("D7").Value = "-1"
$ws.Range("D8").Value = "-1"

# --- Style updates: D2:D6 and D7:D8 move from style s=2 to s=6 (grey-fill text style) ---
# Donor cell D22 already carries style index 6; copy formats only (keeps values intact).
$ws.Range("D22").Copy() | Out-Null
$ws.Range("D2:D8").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# --- Sheet view: active cell / scroll position ---
$ws.Activate()
$ws.Range("G50").Select() | Out-Null
